$d = $word.ActiveDocument
$r = $d.Range(780, 781)
$xml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:tblPr><w:tblStyle w:val="Tabelacomgrade"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4986"/><w:gridCol w:w="4961"/></w:tblGrid><w:tr w:rsidR="00A37AEE" w14:paraId="1A340D7B" w14:textId="77777777" w:rsidTr="00B923B9"><w:tc><w:tcPr><w:tcW w:w="4986" w:type="dxa"/></w:tcPr><w:p w14:paraId="68F07D5B" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"><w:r><w:t>n</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4961" w:type="dxa"/></w:tcPr><w:p w14:paraId="6D05FFD2" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"/></w:tc></w:tr><w:tr w:rsidR="00A37AEE" w14:paraId="36C79F25" w14:textId="77777777" w:rsidTr="00B923B9"><w:tc><w:tcPr><w:tcW w:w="4986" w:type="dxa"/></w:tcPr><w:p w14:paraId="59F5FAA1" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"><w:proofErr w:type="spellStart"/><w:r><w:t>Σx</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4961" w:type="dxa"/></w:tcPr><w:p w14:paraId="4C9CCF40" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"/></w:tc></w:tr><w:tr w:rsidR="00A37AEE" w14:paraId="598DBAAA" w14:textId="77777777" w:rsidTr="00B923B9"><w:tc><w:tcPr><w:tcW w:w="4986" w:type="dxa"/></w:tcPr><w:p w14:paraId="23312091" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"><w:r><w:t>Σx²</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4961" w:type="dxa"/></w:tcPr><w:p w14:paraId="2EBB1126" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"/></w:tc></w:tr><w:tr w:rsidR="00A37AEE" w14:paraId="133E344D" w14:textId="77777777" w:rsidTr="00B923B9"><w:tc><w:tcPr><w:tcW w:w="4986" w:type="dxa"/></w:tcPr><w:p w14:paraId="2D3F001F" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"><w:proofErr w:type="spellStart"/><w:r><w:t>Σy</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4961" w:type="dxa"/></w:tcPr><w:p w14:paraId="64F48278" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"/></w:tc></w:tr><w:tr w:rsidR="00A37AEE" w14:paraId="75996487" w14:textId="77777777" w:rsidTr="00B923B9"><w:tc><w:tcPr><w:tcW w:w="4986" w:type="dxa"/></w:tcPr><w:p w14:paraId="5469275E" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"><w:r><w:t>Σ(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4961" w:type="dxa"/></w:tcPr><w:p w14:paraId="6118203B" w14:textId="77777777" w:rsidR="00A37AEE" w:rsidRDefault="00A37AEE" w:rsidP="00A37AEE"/></w:tc></w:tr></w:tbl>'
$r.InsertXML($xml)
